# Apply "calculated sample 106 concentration after speed vacuum" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace sample 106's placeholder "N/A" cells with real formulas first so the
# shared string table frees up that slot before the new labels are written.
$ws.Range("F10").Formula = "=(F15*F14)/25"
$ws.Range("H10").Formula = "=(G10*F10)/E10"

# New labels used for the speed-vacuum recalculation block (E14:E16)
$ws.Range("E15").Value = "New concentration (ng/µL)"
$ws.Range("E14").Value = "106 Volume after speed vacuum (µL)"
$ws.Range("E16").Value = "Volume water needed to dilute to 25 ng/µL (µL)"

# Column E needs a bit more width for the new labels (stored width snaps to 13)
$ws.Columns.Item(5).ColumnWidth = 12.1

# Style to match: bold header font, wrap text, no explicit horizontal alignment.
# Build the style on E14 only, then format-paint it onto E15:E16 so the
# workbook doesn't accumulate orphaned intermediate cell styles.
$ws.Range("E14").WrapText = $true
$ws.Range("E14").Font.Bold = $true
$ws.Range("E14").Copy()
$ws.Range("E15:E16").PasteSpecial(-4122)
$ws.Range("E14").RowHeight = 48
$ws.Range("E15").RowHeight = 48
$ws.Range("E16").RowHeight = 64

# Values / formulas for the new block
$ws.Range("F14").Value = 18
$ws.Range("F15").Formula = "=(B10*C10)/18"
$ws.Range("F16").Formula = "=F10-F14"

# Update selection / view to match the final state
$ws.Range("E17").Select()
